# Update the first trade row (row 3) with new Quantity and Average Credit values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 3.85
$ws.Range("E3").Value = 3

# Update the active selection to reflect the new cursor position.
$ws.Activate()
$ws.Range("G18").Select()
